{"js": "const replacements = [\n    [\"2023-09-30 Saturday\", \"2023-10-01 Sunday\"],\n    [\"74\u00d796=7104\", \"22\u00d718=396\"],\n    [\"61\u00d788=5368\", \"35\u00d764=2240\"],\n    [\"26\u00d747=1222\", \"31\u00d733=1023\"],\n    [\"92\u00d796=8832\", \"39\u00d756=2184\"],\n    [\"43\u00d777=3311\", \"23\u00d731=713\"],\n    [\"98\u00d713=1274\", \"30\u00d774=2220\"],\n    [\"69\u00d723=1587\", \"70\u00d770=4900\"],\n    [\"44\u00d780=3520\", \"41\u00d754=2214\"],\n    [\"30\u00d780=2400\", \"64\u00d739=2496\"],\n    [\"17\u00d771=1207\", \"85\u00d727=2295\"],\n    [\"58\u00d771=4118\", \"73\u00d760=4380\"],\n    [\"91\u00d769=6279\", \"95\u00d733=3135\"],\n    [\"62\u00d751=3162\", \"88\u00d716=1408\"],\n    [\"15\u00d790=1350\", \"20\u00d779=1580\"],\n    [\"42\u00d754=2268\", \"20\u00d788=1760\"],\n    [\"71\u00d784=5964\", \"73\u00d784=6132\"],\n    [\"86\u00d784=7224\", \"96\u00d779=7584\"],\n    [\"59\u00d780=4720\", \"96\u00d776=7296\"],\n    [\"87\u00d728=2436\", \"30\u00d713=390\"],\n    [\"47\u00d745=2115\", \"89\u00d711=979\"],\n    [\"21\u00d730=630\", \"56\u00d758=3248\"],\n    [\"84\u00d795=7980\", \"56\u00d787=4872\"],\n    [\"14\u00d766=924\", \"68\u00d766=4488\"],\n    [\"55\u00d721=1155\", \"61\u00d778=4758\"],\n    [\"80\u00d717=1360\", \"76\u00d769=5244\"]\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-09-30 Saturday\", \"2023-10-01 Sunday\"),\n    @(\"74\u00d796=7104\", \"22\u00d718=396\"),\n    @(\"61\u00d788=5368\", \"35\u00d764=2240\"),\n    @(\"26\u00d747=1222\", \"31\u00d733=1023\"),\n    @(\"92\u00d796=8832\", \"39\u00d756=2184\"),\n    @(\"43\u00d777=3311\", \"23\u00d731=713\"),\n    @(\"98\u00d713=1274\", \"30\u00d774=2220\"),\n    @(\"69\u00d723=1587\", \"70\u00d770=4900\"),\n    @(\"44\u00d780=3520\", \"41\u00d754=2214\"),\n    @(\"30\u00d780=2400\", \"64\u00d739=2496\"),\n    @(\"17\u00d771=1207\", \"85\u00d727=2295\"),\n    @(\"58\u00d771=4118\", \"73\u00d760=4380\"),\n    @(\"91\u00d769=6279\", \"95\u00d733=3135\"),\n    @(\"62\u00d751=3162\", \"88\u00d716=1408\"),\n    @(\"15\u00d790=1350\", \"20\u00d779=1580\"),\n    @(\"42\u00d754=2268\", \"20\u00d788=1760\"),\n    @(\"71\u00d784=5964\", \"73\u00d784=6132\"),\n    @(\"86\u00d784=7224\", \"96\u00d779=7584\"),\n    @(\"59\u00d780=4720\", \"96\u00d776=7296\"),\n    @(\"87\u00d728=2436\", \"30\u00d713=390\"),\n    @(\"47\u00d745=2115\", \"89\u00d711=979\"),\n    @(\"21\u00d730=630\", \"56\u00d758=3248\"),\n    @(\"84\u00d795=7980\", \"56\u00d787=4872\"),\n    @(\"14\u00d766=924\", \"68\u00d766=4488\"),\n    @(\"55\u00d721=1155\", \"61\u00d778=4758\"),\n    @(\"80\u00d717=1360\", \"76\u00d769=5244\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $r = $d.Content\n    $found = $r.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Search text not found: $oldText\"\n    }\n}\n"}
